$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet/tab (title reflects "through" date)
$ws.Name = "Through 2022-05-06"

# Update the "May (through 05-05)" label to "May (through 05-06)"
$ws.Range("A6").Value = "May (through 05-06)"

# Update May row (row 6) values
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 22
$ws.Range("I6").Value = 19

# Update Total row (row 7) values
$ws.Range("B7").Value = 92
$ws.Range("C7").Value = 171
$ws.Range("D7").Value = 263
$ws.Range("E7").Value = 252
$ws.Range("F7").Value = 164
$ws.Range("G7").Value = 274
$ws.Range("H7").Value = 545
$ws.Range("I7").Value = 570
